# Add an "Address" column between the existing "Name" (E) and "District" (F)
# columns. Inserting a whole column shifts the old F ("District") data into
# the new G column automatically, then we populate the freshly inserted F
# column with the per-row address text extracted from column B/E's combined
# "name + address" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column F ("District"), pushing the
# district data (and everything after it) one column to the right.
$ws.Columns("F").Insert()

# Header
$ws.Range("F2").Value = 'Address'

# Per-row address values.
$ws.Range("F3").Value = 'G H P S No.33 Vantamuri colony'
$ws.Range("F4").Value = 'Govt. High SchoolSangreshkoppaSavadatti'
$ws.Range("F5").Value = 'G H S HirebudanurSavadatti'
$ws.Range("F6").Value = 'S R High School HoskotiRamdurga'
$ws.Range("F7").Value = 'G H S MoogabasavaBail-Hongal'
$ws.Range("F8").Value = 'Govt. High School ChikkadinakoppaKhanapur'
$ws.Range("F9").Value = 'Govt. High School DevalapurBailhongal'
$ws.Range("F10").Value = 'Janata High School SambraBalekundri'
$ws.Range("F11").Value = 'Sharada Girls H S Halaga'
$ws.Range("F12").Value = 'Vidya Mandir High School NesaragiBailhongal'
$ws.Range("F13").Value = 'S N V V S H S Bailhongal'
$ws.Range("F14").Value = 'G H P S HalegattiRamadurga'
$ws.Range("F15").Value = 'Govt. High SchoolRamadurg'
$ws.Range("F16").Value = 'G H S MadlurSavadatti'
$ws.Range("F17").Value = 'G H S NandihalRamdurga'
$ws.Range("F18").Value = 'M H S Yellur'
$ws.Range("F19").Value = 'G H S KatridaddiKulavalliBailhongal'
$ws.Range("F20").Value = 'Ganebail High School GanebailKhanapur'
$ws.Range("F21").Value = 'Shri Dnyaneshwar Vidyalay LokoliKhanapur'
$ws.Range("F22").Value = 'V R M M High SchoolBelavadiBailhongal'
$ws.Range("F23").Value = 'Govt. MLA Model SchoolGurlhosureSaundatti'
# Row 24: source text didn't separate cleanly into name/address, so the
# address cell is left blank (matches upstream data).
$ws.Range("F24").Value = ''
# Row 25: same as row 24 - left blank.
$ws.Range("F25").Value = ''
$ws.Range("F26").Value = 'S B High School TigadolliBailhongal'
$ws.Range("F27").Value = 'Govt. High SchoolVeerapur  Kittur RangeBailhongal'
$ws.Range("F28").Value = 'Govt. High School SalapurRamadurg'
$ws.Range("F29").Value = 'S F S Comp PU College (High School Section) SurebanRamdurg'
$ws.Range("F30").Value = 'Mahantesh High School SalahalliRamdurg'
$ws.Range("F31").Value = 'S M S High SchoolMurgud'
$ws.Range("F32").Value = 'ST. Paul’s High School Camp'
$ws.Range("F33").Value = 'G H S KaghadalSavadatti'
# Row 34: left blank, same reason as rows 24/25.
$ws.Range("F34").Value = ''
$ws.Range("F35").Value = 'Govt. High School HarugoppaSavadatti'
$ws.Range("F36").Value = 'G H S ChikkabellikattiBailhongal'
$ws.Range("F37").Value = 'G H S HoolikottiSavadatti'
$ws.Range("F38").Value = 'K H S Mannar'
$ws.Range("F39").Value = 'Govt. High School MallapurK N Bailhongal'
$ws.Range("F40").Value = 'Govt. Sardar’s PU College'
$ws.Range("F41").Value = 'G H P S ChivatagundiBailwadBaihongal'
$ws.Range("F42").Value = 'Govt High School SunnalRamadurg'
$ws.Range("F43").Value = 'G H S ObalapurRamdurga'
$ws.Range("F44").Value = 'Govt. High School Hudali'
$ws.Range("F45").Value = 'Bashiban High School'
$ws.Range("F46").Value = 'G H S TorangattiRamadurga'
$ws.Range("F47").Value = 'G H S Khanagavi BK Belgaumrural'
$ws.Range("F48").Value = 'M S S H MullurRamdurga'
$ws.Range("F49").Value = 'Shri Shivabasaveshwar'
$ws.Range("F50").Value = 'NaganurBailhongal'
$ws.Range("F51").Value = 'Govt. High School SalahalliRamadurg'
# Row 52: left blank, same reason as rows 24/25/34.
$ws.Range("F52").Value = ''
$ws.Range("F53").Value = 'Govt. High SchoolSulebhavi'
$ws.Range("F54").Value = 'Govt. High School NayanagarBailhongal'
$ws.Range("F55").Value = 'Govt. High SchoolK K Koppa'

# The workbook's recorded used-range extends one column further (through H)
# than any populated cell - touch H55 (a no-op formatting nudge, no value)
# so the saved sheet dimension matches ("A1:H55") without adding real data.
$ws.Range("H55").Font.Bold = $false
